# Fixing matricula of Matc65
# Column A (matricula) values for rows 14-39 need to be corrected.
# Force the range to text format first so the numeric-looking matricula
# IDs stay text (matching the original inlineStr/string cell type)
# instead of being auto-coerced to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14:A39").NumberFormat = "@"

$ws.Range("A14").Value = "217216526"
$ws.Range("A15").Value = "216117974"
$ws.Range("A16").Value = "221117463"
$ws.Range("A17").Value = "217125254"
$ws.Range("A18").Value = "219218129"
$ws.Range("A19").Value = "218215397"
$ws.Range("A20").Value = "220117282"
$ws.Range("A21").Value = "219217429"
$ws.Range("A22").Value = "216216087"
$ws.Range("A23").Value = "220121412"
$ws.Range("A24").Value = "210201260"
$ws.Range("A25").Value = "201520233"
$ws.Range("A26").Value = "217117994"
$ws.Range("A27").Value = "219118481"
$ws.Range("A28").Value = "221119218"
$ws.Range("A29").Value = "219215012"
$ws.Range("A30").Value = "219121541"
$ws.Range("A31").Value = "214007731"
$ws.Range("A32").Value = "219215013"
$ws.Range("A33").Value = "220117290"
$ws.Range("A34").Value = "219118473"
$ws.Range("A35").Value = "220117273"
$ws.Range("A36").Value = "220120071"
$ws.Range("A37").Value = "221216783"
$ws.Range("A38").Value = "214120645"
$ws.Range("A39").Value = "220217140"

# Reset style so cells don't carry a lingering custom number-format index
$ws.Range("A14:A39").Style = "Normal"
